$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "42.415.59"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.310.33"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue "D5" "319.36"
$ws.Range("E5").Value = "  +3.32%  "

Set-TextValue "D6" "103.81"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("E7").Value = "  +1.24%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  +1.68%  "

Set-TextValue "D10" "39.85"
$ws.Range("E10").Value = "  +0.92%  "

Set-TextValue "D11" "0.0911"
$ws.Range("E11").Value = "  +0.44%  "

Set-TextValue "D12" "8.35"
$ws.Range("E12").Value = "  +1.37%  "

Set-TextValue "D13" "0.106"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("E14").Value = "  +1.25%  "

Set-TextValue "D15" "15.39"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "2.660.56"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "2.312.74"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "42.610.39"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("E19").Value = "  +1.28%  "

Set-TextValue "D20" "0.0000105"
$ws.Range("E20").Value = "  +1.28%  "

Set-TextValue "D21" "73.42"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D22" "281.92"
$ws.Range("E22").Value = "  +9.14%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D23" "3.60"
$ws.Range("E23").Value = "  +4.53%  "

Set-TextValue "D24" "10.76"
$ws.Range("E24").Value = "  +17.99%  "

$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("E26").Value = "  -0.25%  "

Set-TextValue "D27" "10.92"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").Value = "  +4.64%  "

Set-TextValue "D29" "22.97"
$ws.Range("E29").Value = "  +0.71%  "

Set-TextValue "D30" "36.18"
$ws.Range("E30").Value = "  +1.89%  "

Set-TextValue "D31" "164.89"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("E32").Value = "  -0.82%  "

Set-TextValue "D33" "5.91"
$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("E34").Value = "  +6.25%  "

Set-TextValue "D35" "2.63"
$ws.Range("E35").Value = "  -9.42%  "

Set-TextValue "D36" "0.116"
$ws.Range("E36").Value = "  -0.04%  "

Set-TextValue "D37" "4.62"
$ws.Range("E37").Value = "  +4.03%  "

Set-TextValue "D38" "0.0365"
$ws.Range("E38").Value = "  +4.95%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "3.71"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D40" "2.77"
$ws.Range("E40").Value = "  +5.86%  "

$ws.Range("E41").Value = "  +3.53%  "

Set-TextValue "D42" "99.98"
$ws.Range("E42").Value = "  +0.20%  "

Set-TextValue "D43" "69.90"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("E45").Value = "  -0.08%  "

Set-TextValue "D46" "12.15"
$ws.Range("E46").Value = "  +1.85%  "

Set-TextValue "D47" "113.16"
$ws.Range("E47").Value = "  +2.34%  "

Set-TextValue "D48" "78.97"
$ws.Range("E48").Value = "  +7.40%  "

Set-TextValue "D49" "8.96"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("D51").Value = "1.619.67"
$ws.Range("E51").Value = "  +6.41%  "
